$d = $word.ActiveDocument

$replacements = @(
    @("84×40=3360", "86×46=3956"),
    @("25×84=2100", "49×30=1470"),
    @("96×88=8448", "46×76=3496"),
    @("67×15=1005", "25×44=1100"),
    @("68×30=2040", "65×61=3965"),
    @("56×38=2128", "76×86=6536"),
    @("72×24=1728", "57×34=1938"),
    @("80×73=5840", "15×60=900"),
    @("80×52=4160", "87×14=1218"),
    @("14×15=210",  "72×62=4464"),
    @("99×92=9108", "15×88=1320"),
    @("12×41=492",  "22×55=1210"),
    @("74×54=3996", "23×85=1955"),
    @("65×34=2210", "13×60=780"),
    @("39×81=3159", "85×34=2890"),
    @("55×12=660",  "25×94=2350"),
    @("42×89=3738", "92×24=2208"),
    @("62×99=6138", "42×64=2688"),
    @("25×63=1575", "11×95=1045"),
    @("93×43=3999", "77×99=7623"),
    @("43×32=1376", "48×72=3456"),
    @("84×62=5208", "30×54=1620"),
    @("47×48=2256", "86×27=2322"),
    @("38×57=2166", "57×66=3762"),
    @("29×40=1160", "83×35=2905")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
